$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Task list - Core library")

# --- Add new task row (row 10 / task #9) ---
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "low"
$ws.Range("C10").Value = "Implement an error handling logic to the whole application"

$ws.Range("D10").Value = "There should be a Protect-like method for handling errors and every exception should be caught. "
$chars = $ws.Range("D10").Characters(19, 7)
$chars.Font.Bold = $true

$ws.Range("E10").Value = "to-do"
$ws.Range("F10").Value = 42059

$ws.Rows.Item(10).RowHeight = 30

# --- Update view: scroll + selection ---
$ws.Range("E11").Select()
$ActiveWindow = $excel.ActiveWindow
$ActiveWindow.ScrollRow = 4
